# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar mensaje de conversión del día ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.5 = 30400.28 pesos`n✅ 30400.28 pesos = 7.48 = 968.01 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas: actualizar tasas N10, N12 y O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 133.408
$ws2.Range("N12").Value = 4063.79
$ws2.Range("O12").Value = 129.4
